$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.994.80"
$ws.Range("E2").Value = "  -4.84%  "
$ws.Range("D3").Value = "2.209.44"
$ws.Range("E3").Value = "  -7.61%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "295.74"
$ws.Range("C5").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = "  -5.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "80.13"
$ws.Range("C6").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E6").Value = "  -10.00%  "
$ws.Range("E7").Value = "  -5.48%  "
$ws.Range("E9").Value = "  -8.40%  "
$ws.Range("E10").Value = "  -8.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "27.90"
$ws.Range("C11").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("E11").Value = "  -11.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.74"
$ws.Range("C12").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("E12").Value = "  -13.12%  "
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("D14").Value = "2.556.37"
$ws.Range("E14").Value = "  -7.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.08"
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "  -8.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.84"
$ws.Range("C16").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "  -9.43%  "
$ws.Range("D17").Value = "2.233.21"
$ws.Range("E17").Value = "  -7.17%  "
$ws.Range("E18").Value = "  -8.26%  "
$ws.Range("D19").Value = "38.875.96"
$ws.Range("E19").Value = "  -5.15%  "
$ws.Range("E20").Value = "  -6.62%  "
$ws.Range("E21").Value = "  -8.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.57"
$ws.Range("C22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "  -7.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.76"
$ws.Range("C23").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = "  -10.43%  "
$ws.Range("E24").Value = "  -4.17%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.38"
$ws.Range("C26").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = "  -10.73%  "
$ws.Range("E27").Value = "  -5.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.20"
$ws.Range("C28").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "  -7.63%  "
$ws.Range("E29").Value = "  -2.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.86"
$ws.Range("C30").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "  -5.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "148.10"
$ws.Range("C31").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "  -3.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.00"
$ws.Range("C32").Copy()
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("E32").Value = "  -9.51%  "
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("E34").Value = "  -9.88%  "
$ws.Range("E35").Value = "  -4.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0685"
$ws.Range("C36").Copy()
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("E36").Value = "  -7.25%  "
$ws.Range("E37").Value = "  -5.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0952"
$ws.Range("C38").Copy()
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("E38").Value = "  -5.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.61"
$ws.Range("C39").Copy()
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("E39").Value = "  -7.86%  "
$ws.Range("E40").Value = "  -9.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.24"
$ws.Range("C41").Copy()
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("E41").Value = "  -11.80%  "
$ws.Range("E42").Value = "  -7.25%  "
$ws.Range("D43").Value = "1.891.80"
$ws.Range("E43").Value = "  -4.30%  "
$ws.Range("E44").Value = "  -12.01%  "
$ws.Range("E45").Value = "  -7.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.98"
$ws.Range("C46").Copy()
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("E46").Value = "  -10.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.79"
$ws.Range("C47").Copy()
$ws.Range("D47").PasteSpecial(-4122)
$ws.Range("E47").Value = "  -9.38%  "
$ws.Range("E48").Value = "  -10.49%  "
$ws.Range("D49").Value = "2.425.61"
$ws.Range("E49").Value = "  -7.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.77"
$ws.Range("C50").Copy()
$ws.Range("D50").PasteSpecial(-4122)
$ws.Range("E50").Value = "  -4.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.08"
$ws.Range("C51").Copy()
$ws.Range("D51").PasteSpecial(-4122)
$ws.Range("E51").Value = "  -0.71%  "
$excel.CutCopyMode = $false
